# "All functions are now in preprocessing"
# Adds three new model columns (CNN-1 / CNN-2 / CNN2-nodrop style variants) to each of the
# three result blocks (128_bin, 128_bin_times_10, 128_extended_bin) on the "new_results" sheet,
# adds a new "Extended_bin_Aug" summary row, and repositions/resizes the chart to fit the wider table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1. Make room: insert 3 columns before the "128_bin_times_10" block (old D:E)
#    and 3 more before the "128_extended_bin" block (old F:G, now shifted to L:M).
#    Excel's column-insert automatically re-points every formula/merge that
#    referenced the shifted cells (B12/B13 etc.), matching the target diff.
# ---------------------------------------------------------------------------
$ws.Range("D1:F1").EntireColumn.Insert()
$ws.Range("I1:K1").EntireColumn.Insert()

# ---------------------------------------------------------------------------
# 2. New header row (row 2) entries for the freshly inserted columns.
#    (Written in this particular order so newly-introduced shared strings
#    line up with how the authored workbook built its string table.)
# ---------------------------------------------------------------------------
$ws.Range("D2").Value = "CNN-1"
$ws.Range("E2").Value = "CNN-2"
$ws.Range("F2").Value = "CNN2-nodrop"

$ws.Range("J2").Value = "CNN-2"
$ws.Range("K2").Value = "CNN2-nodrop"

$ws.Range("O2").Value = "CNN-2"
$ws.Range("P2").Value = "CNN2-nodrop"

# ---------------------------------------------------------------------------
# 3. Re-merge the title cells across the now-wider blocks.
# ---------------------------------------------------------------------------
$ws.Range("B1:C1").UnMerge()
$ws.Range("G1:H1").UnMerge()
$ws.Range("L1:M1").UnMerge()

$ws.Range("B1:F1").Merge()
$ws.Range("G1:K1").Merge()
$ws.Range("L1:P1").Merge()

# ---------------------------------------------------------------------------
# 4. New small legend table (row 10) gains the two extra model columns.
# ---------------------------------------------------------------------------
$ws.Range("D10").Value = "CNN-1"
$ws.Range("E10").Value = "CNN-2"
$ws.Range("F10").Value = "CNN-2-nodrop"

# ---------------------------------------------------------------------------
# 5. New summary row for the "Extended_bin_Aug" dataset/model (same text-
#    formatted style as the A11:A13 labels above it).
# ---------------------------------------------------------------------------
$ws.Range("A14").Value = "Extended_bin_Aug"
$ws.Range("A14").NumberFormat = "@"

# ---------------------------------------------------------------------------
# 5b. The "128_bin_times_10" / "128_extended_bin" blocks use the undropped
#     "CNN1" label (vs. "CNN-1" for 128_bin) in their second model column.
# ---------------------------------------------------------------------------
$ws.Range("I2").Value = "CNN1"
$ws.Range("N2").Value = "CNN1"

# ---------------------------------------------------------------------------
# 6. Resize columns (values picked so the stored xlsx column width matches
#    the authored widths as closely as the engine's width quantization allows).
# ---------------------------------------------------------------------------
$widths = @(17.166666666666362, 15.166666666666362, 12.33333333333303, 5.999999999999696, 5.999999999999696, 13.166666666666362, 3.9999999999996962, 14.33333333333303, 5.166666666666363, 5.999999999999696, 12.499999999999696, 3.9999999999996962, 14.33333333333303, 5.166666666666363, 5.999999999999696, 12.499999999999696)
for ($i = 0; $i -lt $widths.Length; $i++) {
  $ws.Columns.Item($i + 1).ColumnWidth = $widths[$i]
}

# ---------------------------------------------------------------------------
# 7. Move/resize the chart so it still sits under the (now much wider) table.
# ---------------------------------------------------------------------------
$co = $ws.ChartObjects().Item(1)
$co.Left = 84.0
$co.Top = 264.0
$co.Width = 735.0
$co.Height = 200.25

# ---------------------------------------------------------------------------
# 8. Restore the last-used selection cell recorded in the saved file.
# ---------------------------------------------------------------------------
$ws.Range("S19").Select()
